# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the combined "全部类型" sheet to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 7763
$wsExhibit.Range("F6").Value  = 205
$wsExhibit.Range("F10").Value = 446
$wsExhibit.Range("F13").Value = 440
$wsExhibit.Range("F15").Value = 66
$wsExhibit.Range("F17").Value = 5626
$wsExhibit.Range("F19").Value = 220
$wsExhibit.Range("F20").Value = 1039
$wsExhibit.Range("F21").Value = 229
$wsExhibit.Range("F22").Value = 327

# Sheet "全部类型" (All types) - mirrors the exhibition rows plus the
# interleaved performance rows, so row numbers shift by one after row 16.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 7763
$wsAll.Range("F6").Value  = 205
$wsAll.Range("F10").Value = 446
$wsAll.Range("F13").Value = 440
$wsAll.Range("F15").Value = 66
$wsAll.Range("F18").Value = 5626
$wsAll.Range("F21").Value = 220
$wsAll.Range("F22").Value = 1039
$wsAll.Range("F23").Value = 229
$wsAll.Range("F24").Value = 327
